# edit.ps1 - apply the tracked changes to assignment.docx
#
# 1. Remove the stray "_GoBack" bookmark from the top of the document.
# 2. "Contrapositive form of 3rd statement" -> "Contrapositive form of 1st statement"
#    and re-drop the "_GoBack" bookmark right after the ordinal suffix run
#    (this is where Word's last-edit-position bookmark ends up).
# 3. Split "If a bird is an Ostrich, then it does not live on mince pies"
#    into "If a bird is an Ostrich, then it " + "is at least 9 feet tall",
#    keeping both halves as separate runs with identical run formatting.

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark wherever it currently lives ---------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. "Contrapositive form of 3" / "rd" -> "... of 1" / "st" ------------
$d.Content.Find.Execute("Contrapositive form of 3", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Contrapositive form of 1", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Contrapositive form of 1", $true, $false, $false, $false,
                   $false, $true, 1, $false, "", 0) | Out-Null
$suffixStart = $rng.End
$suffixRng = $d.Range($suffixStart, $suffixStart + 2)
$suffixRng.Text = "st"

# Re-insert _GoBack immediately after the new "st" run (collapsed range).
$goBackPoint = $d.Range($suffixStart + 2, $suffixStart + 2)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# --- 3. Split the Ostrich sentence into two runs ---------------------------
$tail = $d.Content
$tail.Find.Execute("does not live on mince pies", $true, $false, $false,
                    $false, $false, $true, 1, $false, "", 0) | Out-Null

# Toggling a character property around the text replacement keeps this
# edit from being re-merged into the preceding run, even though the final
# formatting (no bold) ends up identical on both sides.
$tail.Font.Bold = 1
$tail.Text = "is at least 9 feet tall"
$tail.Font.Bold = 0
